$wb = $excel.ActiveWorkbook

# Regular_Timetable
$ws = $wb.Worksheets.Item("Regular_Timetable")
$ws.Range("B3").Value = "MA261 [C001]"
$ws.Range("C3").Value = "MA261 [C001]"
$ws.Range("D3").Value = "EC261 [C104]"
$ws.Range("E3").Value = "EC262 [C304]"
$ws.Range("F3").Value = "EC263 [C305]"
$ws.Range("B4").Value = "MA262 [C001]"
$ws.Range("C4").Value = "MA262 [C001]"
$ws.Range("D4").Value = "EC262 [C304]"
$ws.Range("E4").Value = "EC263 [C305]"
$ws.Range("F4").Value = "CS307 [C302]"
$ws.Range("C6").Value = "MA263 [C101]"
$ws.Range("E6").Value = "EC263 (Lab) [L207]"
$ws.Range("B7").Value = "MA261 (Tutorial) [C001]"
$ws.Range("E7").Value = "EC263 (Lab) [L207]"
$ws.Range("B8").Value = "MA263 [C101]"
$ws.Range("C8").Value = "EC261 [C204]"
$ws.Range("D8").Value = "EC262 (Lab) [L207]"
$ws.Range("E8").Value = "CS307 [C302]"
$ws.Range("B9").Value = "MA262 (Tutorial) [C001]"
$ws.Range("D9").Value = "EC262 (Lab) [L207]"
$ws.Range("D26").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D27").Value = "Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D32").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E32").Value = "Tue 14:30-15:30 [C205]"

# PreMid_Timetable
$ws = $wb.Worksheets.Item("PreMid_Timetable")
$ws.Range("B3").Value = "EC261 [C204]"
$ws.Range("C3").Value = "EC261 [C204]"
$ws.Range("D3").Value = "MA263 [C101]"
$ws.Range("E3").Value = "MA263 [C104]"
$ws.Range("B4").Value = "CS307 [C302]"
$ws.Range("C4").Value = "CS307 [C302]"
$ws.Range("D4").Value = "EC262 [C305]"
$ws.Range("E4").Value = "EC262 [C102]"
$ws.Range("B7").Value = "EC261 (Tutorial) [C104]"
$ws.Range("C7").Value = "MA261 (Tutorial) [C001]"
$ws.Range("E8").Value = "EC262 (Lab) [L207]"
$ws.Range("E9").Value = "EC262 (Lab) [L207]"
$ws.Range("D26").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D27").Value = "Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D32").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E32").Value = "Tue 14:30-15:30 [C205]"

# PostMid_Timetable
$ws = $wb.Worksheets.Item("PostMid_Timetable")
$ws.Range("B3").Value = "EC261 [C203]"
$ws.Range("C3").Value = "EC261 [C203]"
$ws.Range("D3").Value = "EC263 [C001]"
$ws.Range("E3").Value = "EC263 [C305]"
$ws.Range("B4").Value = "CS307 [C205]"
$ws.Range("C4").Value = "CS307 [C203]"
$ws.Range("B6").Value = "MA262 [C001]"
$ws.Range("C7").Value = "MA262 (Tutorial) [C001]"
$ws.Range("C8").Value = "MA262 [C001]"
$ws.Range("B9").Value = "CS307 (Tutorial) [C203]"
$ws.Range("D26").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D27").Value = "Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D32").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E32").Value = "Tue 14:30-15:30 [C205]"
